$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper to write a full "match row" (columns B..AC) from a hashtable of
# column-letter -> value. Columns not present in the hashtable are left
# untouched (so callers only need to specify what actually changes for the
# "update odds / final score" rows, or everything for the "swap" rows).
# ---------------------------------------------------------------------------
function Set-RowValues {
    param(
        [int]$Row,
        [hashtable]$Values
    )
    foreach ($col in $Values.Keys) {
        $ws.Range("$col$Row").Value = $Values[$col]
    }
}

# ---------------------------------------------------------------------------
# Rows 95 / 96 : the two match records were swapped (everything except the
# running index in column A).
# ---------------------------------------------------------------------------
$cols = @('B','C','D','E','F','G','H','I','J','K','L','M','N','O','P','Q','R','S','T','U','V','W','X','Y','Z','AA','AB','AC')

function Swap-Rows {
    param([int]$RowA, [int]$RowB)
    $valsA = @{}
    $valsB = @{}
    foreach ($c in $cols) {
        $valsA[$c] = $ws.Range("$c$RowA").Value()
        $valsB[$c] = $ws.Range("$c$RowB").Value()
    }
    foreach ($c in $cols) {
        $ws.Range("$c$RowA").Value = $valsB[$c]
        $ws.Range("$c$RowB").Value = $valsA[$c]
    }
}

Swap-Rows -RowA 95 -RowB 96

# ---------------------------------------------------------------------------
# Rows 129 / 130 / 131 : 3-way rotation -> new129 = old131, new130 = old129,
# new131 = old130.
# ---------------------------------------------------------------------------
$v129 = @{}
$v130 = @{}
$v131 = @{}
foreach ($c in $cols) {
    $v129[$c] = $ws.Range("${c}129").Value()
    $v130[$c] = $ws.Range("${c}130").Value()
    $v131[$c] = $ws.Range("${c}131").Value()
}
foreach ($c in $cols) {
    $ws.Range("${c}129").Value = $v131[$c]
    $ws.Range("${c}130").Value = $v129[$c]
    $ws.Range("${c}131").Value = $v130[$c]
}

# ---------------------------------------------------------------------------
# Rows 192 / 193 : swapped (like 95/96).
# ---------------------------------------------------------------------------
Swap-Rows -RowA 192 -RowB 193

# ---------------------------------------------------------------------------
# Rows 200 / 201 : swapped (like 95/96).
# ---------------------------------------------------------------------------
Swap-Rows -RowA 200 -RowB 201

# ---------------------------------------------------------------------------
# Rows 224-227 : these matches were played since the last update, so the
# final score (H/I/J) and the closing odds (K..AC) now have real values.
# ---------------------------------------------------------------------------
Set-RowValues -Row 224 -Values @{
    H='0'; I='0'; J='D'
    Q='0.75'; R='1.8'; S='2'; U='1.925'; V='1.875'
    W='-1'; X='2.6'; Y='-1'; Z='0.8'; AA='-1'
    AB='-1'; AC='0.875'
    P='1.75'
}

Set-RowValues -Row 225 -Values @{
    H='3'; I='0'; J='H'
    N='3.4'; O='3.2'; P='2.15'
    U='2'; V='1.8'
    W='2.4'; X='-1'; Y='-1'; Z='0.95'; AA='-1'
    AB='1'; AC='-1'
}

Set-RowValues -Row 226 -Values @{
    H='3'; I='1'; J='H'
    N='1.3'; O='5'; P='7'; Q='-1.5'
    R='2'; S='1.8'
    U='1.975'; V='1.825'
    W='0.3'; X='-1'; Y='-1'; Z='1'; AA='-1'
    AB='0.9750000000000001'; AC='-1'
}

Set-RowValues -Row 227 -Values @{
    H='3'; I='3'; J='D'
    O='4.5'; P='7'
    U='1.85'; V='1.95'
    W='-1'; X='3.5'; Y='-1'; Z='-1'
    AB='0.8500000000000001'; AC='-1'
}
# AA227 is not populated for this match; remove the stray value that was
# there before.
$ws.Range("AA227").ClearContents()

# ---------------------------------------------------------------------------
# Rows 228 / 229 : no longer part of the sheet -> delete them (this also
# shrinks the used range from AC229 down to AC227, as in the target file).
# ---------------------------------------------------------------------------
$ws.Rows.Item(228).Delete()
$ws.Rows.Item(228).Delete()
